$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A22").Value = "2024-11-02 18:08:29"
$ws.Range("B22").Value = "Success"

$ws.Range("A23").Value = "2024-11-02 18:10:15"
$ws.Range("B23").Value = "Success"
